$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Examples")

# Add the new "Tags" column (column O) header, copying the format/style
# from the neighboring N1 cell (same header row style).
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("O1").Value = "Tags"

# Populate the new column's data rows. Write the leading-space variant
# first so it claims the earlier shared-string table slot, matching the
# order strings were interned in the source workbook.
$ws.Range("O3").Value = " @examples"
$ws.Range("O2").Value = "@examples"
$ws.Range("O4").Value = "@examples"
$ws.Range("O5").Value = "@examples"
$ws.Range("O6").Value = "@examples"

# Clear clipboard marching ants / selection artifacts and set the new
# active selection to O6 as in the edited workbook.
$excel.CutCopyMode = 0
$ws.Range("O6").Select()
